$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.014926894567906857
$ws.Range("C2").Value = 0.006258303299546242
$ws.Range("D2").Value = 0.004506985656917095
$ws.Range("E2").Value = 0.003292342182248831
$ws.Range("F2").Value = 0.000024721293812035583
$ws.Range("J2").Value = 0.12727078795433044
$ws.Range("K2").Value = 1.4152871370315552
